$d = $word.ActiveDocument

# --- Step 1: turn the single paragraph into two.
#     Insert a paragraph break right at the end of the document content -
#     that point sits right before the "_GoBack" bookmark, so splitting
#     there moves the (empty) bookmark onto the new second paragraph,
#     exactly like Word does when you press Enter at the end of a line.
$endRng = $d.Content
$endRng.Collapse(0)
$endRng.InsertBefore("`r")

# Type the new line's text at the very start of the (now second) paragraph,
# i.e. before the bookmark that now lives there.
$p2 = $d.Paragraphs.Item(2)
$p2Start = $p2.Range
$p2Start.Collapse(1)
$p2Start.InsertBefore("Outra linha ")

# --- Step 2: on paragraph 1, capitalize the first letter and split the
#     run in two ("L" / rest), bracketed by proofErr spell-check markers -
#     the same artifact Word leaves behind after autocorrecting/checking
#     a word like "lknlknlkn".
$p1 = $d.Paragraphs.Item(1)
$full = $p1.Range
$full.MoveEnd(1, -1) | Out-Null   # exclude the trailing paragraph mark
$original = $full.Text
$firstChar = $original.Substring(0, 1).ToUpper()
$rest = $original.Substring(1)
$full.Delete()

$insPt = $d.Paragraphs.Item(1).Range
$insPt.Collapse(1)
$xml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>$firstChar</w:t></w:r>
<w:r><w:t>$rest</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
$insPt.InsertXML($xml)

Write-Output "Final text: [$($d.Content.Text)]"
